# Row 22 and Row 23 in the "Artfynd" sheet swap their species/observation
# data (columns A, B, D, E, F, G, H, I, Q, R) while the location/date/
# reporter columns stay put (they were already identical between the two
# rows). Re-create that by capturing both rows' old values first, then
# writing each row's new values from the other row's captured data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 22
$row2 = 23

# Columns whose content is exchanged between the two rows.
$cols = "A", "B", "D", "E", "F", "G", "H", "Q", "R"

# Snapshot current values before overwriting anything.
$old1 = @{}
$old2 = @{}
foreach ($col in $cols) {
    $old1[$col] = $ws.Range("$col$row1").Value2
    $old2[$col] = $ws.Range("$col$row2").Value2
}
$oldI1 = $ws.Range("I$row1").Value2
$oldI2 = $ws.Range("I$row2").Value2

# Write the numeric / plain-text columns (simple value swap).
foreach ($col in $cols) {
    $ws.Range("$col$row1").Value = $old2[$col]
    $ws.Range("$col$row2").Value = $old1[$col]
}

# Column I ("Antal") is stored as text even though it holds digit strings
# (e.g. "1"), and row 22 previously held no value there at all. Force the
# destination cells to text so a value like "1" isn't re-interpreted as a
# number, then drop the formatting override so no stray style is left
# behind once the text is in place.
$ws.Range("I$row1").NumberFormat = "@"
$ws.Range("I$row1").Value = $oldI2
$ws.Range("I$row1").ClearFormats()

$ws.Range("I$row2").NumberFormat = "@"
$ws.Range("I$row2").Value = $oldI1
$ws.Range("I$row2").ClearFormats()
